$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 218 (shifts old rows 218-313 down to 219-314),
# inheriting formatting from the row above (row 217), which already
# carries the date style (s="2") on column D.
$ws.Rows.Item(218).Insert()

# Populate the newly inserted row 218 with the new record.
$ws.Cells.Item(218, 1).Value  = 4
$ws.Cells.Item(218, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(218, 3).Value  = "Los Lagos"
$ws.Cells.Item(218, 4).Value  = 44806
$ws.Cells.Item(218, 5).Value  = 10
$ws.Cells.Item(218, 6).Value  = 100112017
$ws.Cells.Item(218, 7).Value  = "Apio"
$ws.Cells.Item(218, 8).Value  = "Americana (o)"
$ws.Cells.Item(218, 9).Value  = "Primera"
$ws.Cells.Item(218, 10).Value = 50
$ws.Cells.Item(218, 11).Value = 14000
$ws.Cells.Item(218, 12).Value = 15000
$ws.Cells.Item(218, 13).Value = 14500
$ws.Cells.Item(218, 14).Value = "$/docena de matas"
$ws.Cells.Item(218, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(218, 16).Value = 2417
$ws.Cells.Item(218, 17).Value = 6
$ws.Cells.Item(218, 18).Value = "Hortaliza"
